$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: reduce numeric precision to 2 decimal places (custom accuracy)
$ws.Range("B5").Value = 15.31
$ws.Range("C5").Value = 11.28
$ws.Range("D5").Value = 1.06
$ws.Range("E5").Value = 33.48
$ws.Range("F5").Value = 27.03
$ws.Range("G5").Value = 11.98
$ws.Range("H5").Value = 47.02
$ws.Range("I5").Value = 18.62
$ws.Range("J5").Value = 8.199999999999999
$ws.Range("K5").Value = 12.06
$ws.Range("L5").Value = 13.4
$ws.Range("M5").Value = 14.25
$ws.Range("N5").Value = 3.86
$ws.Range("O5").Value = 12.03
$ws.Range("P5").Value = 17.06
$ws.Range("Q5").Value = 10.25
$ws.Range("R5").Value = 0.77
$ws.Range("S5").Value = 0.68
$ws.Range("T5").Value = 175.8
$ws.Range("U5").Value = 33.65
$ws.Range("V5").Value = 11.1
$ws.Range("W5").Value = 22.49
$ws.Range("X5").Value = 11.79
$ws.Range("Y5").Value = 1.89
$ws.Range("Z5").Value = 22.83
$ws.Range("AA5").Value = 9.81
$ws.Range("AB5").Value = 8.75
$ws.Range("AC5").Value = 10.29
$ws.Range("AD5").Value = 13.97
$ws.Range("AE5").Value = 0.54
$ws.Range("AF5").Value = 42.83
$ws.Range("AG5").Value = 6.21
$ws.Range("AH5").Value = 13.88

# Narrow columns Q (17) and AC (29) by one character unit
$ws.Columns.Item(17).ColumnWidth = 6.17
$ws.Columns.Item(29).ColumnWidth = 6.17

# Drop the last data row (row 6) entirely
$ws.Rows.Item(6).Delete()
